# Codon optimization workbook: row 7 previously duplicated row 6 (the per-species
# median). The data-selection logic was improved so row 7 now holds its own set
# of distinct per-column figures (counts used by the new selection GUI) instead
# of repeating row 6. Sheet1 cells already carry a Text ("@") number format, so
# assigning numeric-looking strings keeps them stored as text, matching the rest
# of the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A7").Value = "12321"
$ws.Range("B7").Value = "86"
$ws.Range("C7").Value = "49"
$ws.Range("D7").Value = "59"
$ws.Range("E7").Value = "84"
$ws.Range("F7").Value = "51"
$ws.Range("G7").Value = "26"
$ws.Range("H7").Value = "27"
$ws.Range("I7").Value = "45"
$ws.Range("J7").Value = "83"
$ws.Range("K7").Value = "56"
$ws.Range("L7").Value = "30"
$ws.Range("M7").Value = "66"
$ws.Range("N7").Value = "74"
$ws.Range("O7").Value = "30"
$ws.Range("P7").Value = "41"
$ws.Range("Q7").Value = "79"
$ws.Range("R7").Value = "62"
$ws.Range("S7").Value = "49"
$ws.Range("T7").Value = "4"
$ws.Range("U7").Value = "1"
$ws.Range("V7").Value = "43"
$ws.Range("W7").Value = "21"
$ws.Range("X7").Value = "57"
$ws.Range("Y7").Value = "44"
$ws.Range("Z7").Value = "66"
$ws.Range("AA7").Value = "45"
$ws.Range("AB7").Value = "104"
$ws.Range("AC7").Value = "108"
$ws.Range("AD7").Value = "129"
$ws.Range("AE7").Value = "62"
$ws.Range("AF7").Value = "122"
$ws.Range("AG7").Value = "67"
$ws.Range("AH7").Value = "77"
$ws.Range("AI7").Value = "31"
$ws.Range("AJ7").Value = "58"
$ws.Range("AK7").Value = "17"
$ws.Range("AL7").Value = "47"
$ws.Range("AM7").Value = "23"
$ws.Range("AN7").Value = "34"
$ws.Range("AO7").Value = "33"
$ws.Range("AP7").Value = "72"
$ws.Range("AQ7").Value = "23"
$ws.Range("AR7").Value = "38"
$ws.Range("AS7").Value = "32"
$ws.Range("AT7").Value = "72"
$ws.Range("AU7").Value = "48"
$ws.Range("AV7").Value = "35"
$ws.Range("AW7").Value = "45"
$ws.Range("AX7").Value = "49"
$ws.Range("AY7").Value = "26"
$ws.Range("AZ7").Value = "2"
$ws.Range("BA7").Value = "28"
$ws.Range("BB7").Value = "41"
$ws.Range("BC7").Value = "9"
$ws.Range("BD7").Value = "33"
$ws.Range("BE7").Value = "15"
$ws.Range("BF7").Value = "44"
$ws.Range("BG7").Value = "30"
$ws.Range("BH7").Value = "56"
$ws.Range("BI7").Value = "25"
$ws.Range("BJ7").Value = "74"
$ws.Range("BK7").Value = "22"
$ws.Range("BL7").Value = "54"
$ws.Range("BM7").Value = "19"
